$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 87, shifting existing rows 87-161 down to 88-162.
$ws.Rows("87:87").Insert()

# Populate the newly inserted row 87 with the new record's data.
$ws.Cells.Item(87, 1).Value = 11
$ws.Cells.Item(87, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(87, 3).Value = "Bíobío"
$ws.Cells.Item(87, 4).Value = 44729
$ws.Cells.Item(87, 5).Value = 8
$ws.Cells.Item(87, 6).Value = "Fruta"
$ws.Cells.Item(87, 7).Value = 100109
$ws.Cells.Item(87, 8).Value = "Uva"
$ws.Cells.Item(87, 9).Value = 100109001
$ws.Cells.Item(87, 10).Value = "Uva"
$ws.Cells.Item(87, 11).Value = "Red Globe"
$ws.Cells.Item(87, 12).Value = "Primera"
$ws.Cells.Item(87, 13).Value = 100
$ws.Cells.Item(87, 14).Value = 6000
$ws.Cells.Item(87, 15).Value = 7000
$ws.Cells.Item(87, 16).Value = 6500
$ws.Cells.Item(87, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(87, 18).Value = "Región Metropolitana"
$ws.Cells.Item(87, 19).Value = 650
$ws.Cells.Item(87, 20).Value = 10
